# Updates the cryptos list (Sheet1, columns B-E, rows 2-51) with refreshed
# price/volume data, matching the "Updated cryptos list" GitHub Actions commit.
# Two coin pairs also swapped rank order (rows 19/20: Avalanche<->TRON,
# rows 48/49: BabyDogeCoin<->EOS), so Coin/Link/Price/Volume are rewritten there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.025.56"
$ws.Range("E2").Value = "  -1.96%  "
# Row 3
$ws.Range("D3").Value = "2.101.09"
$ws.Range("E3").Value = "  -0.61%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -0.27%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.51"
$ws.Range("E5").Value = "  +1.67%  "
# Row 6
$ws.Range("E6").Value = "  -0.31%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5183"
$ws.Range("E7").Value = "  -1.36%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4483"
$ws.Range("E8").Value = "  -0.38%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09493"
$ws.Range("E9").Value = "  +5.06%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.17"
$ws.Range("E10").Value = "  -2.70%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.175"
$ws.Range("E11").Value = "  +0.52%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.19"
$ws.Range("E12").Value = "  +3.44%  "
# Row 13
$ws.Range("D13").Value = "2.110.72"
$ws.Range("E13").Value = "  -0.62%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.742"
$ws.Range("E14").Value = "  -0.59%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.067"
$ws.Range("E15").Value = "  +0.12%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.36"
$ws.Range("E16").Value = "  +1.50%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001169"
$ws.Range("E17").Value = "  +0.56%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  -0.31%  "
# Row 19
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06710"
$ws.Range("E19").Value = "  +0.19%  "
# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.55"
$ws.Range("E20").Value = "  +6.23%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -0.32%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.184"
$ws.Range("E22").Value = "  -2.70%  "
# Row 23
$ws.Range("D23").Value = "30.117.00"
$ws.Range("E23").Value = "  -1.99%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.67"
$ws.Range("E24").Value = "  -1.11%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  -2.65%  "
# Row 26
$ws.Range("D26").Value = "2.356.70"
$ws.Range("E26").Value = "  -0.64%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.01"
$ws.Range("E27").Value = "  -1.69%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.30"
$ws.Range("E28").Value = "  -0.59%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.529"
$ws.Range("E29").Value = "  -0.51%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.04"
$ws.Range("E30").Value = "  -0.54%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.156"
$ws.Range("E31").Value = "  -3.17%  "
# Row 32
$ws.Range("E32").Value = "  -1.56%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.622"
$ws.Range("E33").Value = "  -0.58%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.250"
$ws.Range("E34").Value = "  -1.93%  "
# Row 35
$ws.Range("E35").Value = "  +0.29%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.163"
$ws.Range("E36").Value = "  +4.65%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.11"
$ws.Range("E37").Value = "  -2.03%  "
# Row 38
$ws.Range("E38").Value = "  -2.89%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06776"
$ws.Range("E39").Value = "  -0.77%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2278"
$ws.Range("E40").Value = "  -1.74%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6938"
$ws.Range("E41").Value = "  +0.95%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.43"
$ws.Range("E42").Value = "  -1.24%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.312"
$ws.Range("E43").Value = "  +4.15%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6687"
$ws.Range("E44").Value = "  +4.13%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.25"
$ws.Range("E45").Value = "  -5.12%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.275"
$ws.Range("E46").Value = "  -1.62%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.639"
$ws.Range("E47").Value = "  -1.73%  "
# Row 48
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.222"
$ws.Range("E48").Value = "  -2.52%  "
# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000342"
$ws.Range("E49").Value = "  -7.29%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.90"
$ws.Range("E50").Value = "  -1.22%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07167"
$ws.Range("E51").Value = "  -1.93%  "
